$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("PIR")
$ws.Range("A312:F323").NumberFormat = "@"
$ws.Cells.Item(312, 1).Value = "2026-02-04"
$ws.Cells.Item(312, 2).Value = "14:27:00"
$ws.Cells.Item(312, 3).Value = "14:00"
$ws.Cells.Item(312, 4).Value = "Bathroom"
$ws.Cells.Item(312, 5).Value = "No Motion"
$ws.Cells.Item(312, 6).Value = "Inactive"
$ws.Cells.Item(313, 1).Value = "2026-02-04"
$ws.Cells.Item(313, 2).Value = "14:27:03"
$ws.Cells.Item(313, 3).Value = "14:00"
$ws.Cells.Item(313, 4).Value = "Bathroom"
$ws.Cells.Item(313, 5).Value = "Motion Detected"
$ws.Cells.Item(313, 6).Value = "Active"
$ws.Cells.Item(314, 1).Value = "2026-02-04"
$ws.Cells.Item(314, 2).Value = "14:27:11"
$ws.Cells.Item(314, 3).Value = "14:00"
$ws.Cells.Item(314, 4).Value = "Bathroom"
$ws.Cells.Item(314, 5).Value = "No Motion"
$ws.Cells.Item(314, 6).Value = "Inactive"
$ws.Cells.Item(315, 1).Value = "2026-02-04"
$ws.Cells.Item(315, 2).Value = "14:27:17"
$ws.Cells.Item(315, 3).Value = "14:00"
$ws.Cells.Item(315, 4).Value = "Bathroom"
$ws.Cells.Item(315, 5).Value = "No Motion"
$ws.Cells.Item(315, 6).Value = "Inactive"
$ws.Cells.Item(316, 1).Value = "2026-02-04"
$ws.Cells.Item(316, 2).Value = "14:27:22"
$ws.Cells.Item(316, 3).Value = "14:00"
$ws.Cells.Item(316, 4).Value = "Bathroom"
$ws.Cells.Item(316, 5).Value = "No Motion"
$ws.Cells.Item(316, 6).Value = "Inactive"
$ws.Cells.Item(317, 1).Value = "2026-02-04"
$ws.Cells.Item(317, 2).Value = "14:27:24"
$ws.Cells.Item(317, 3).Value = "14:00"
$ws.Cells.Item(317, 4).Value = "Bathroom"
$ws.Cells.Item(317, 5).Value = "Motion Detected"
$ws.Cells.Item(317, 6).Value = "Active"
$ws.Cells.Item(318, 1).Value = "2026-02-04"
$ws.Cells.Item(318, 2).Value = "14:27:33"
$ws.Cells.Item(318, 3).Value = "14:00"
$ws.Cells.Item(318, 4).Value = "Bathroom"
$ws.Cells.Item(318, 5).Value = "No Motion"
$ws.Cells.Item(318, 6).Value = "Inactive"
$ws.Cells.Item(319, 1).Value = "2026-02-04"
$ws.Cells.Item(319, 2).Value = "14:27:37"
$ws.Cells.Item(319, 3).Value = "14:00"
$ws.Cells.Item(319, 4).Value = "Bathroom"
$ws.Cells.Item(319, 5).Value = "No Motion"
$ws.Cells.Item(319, 6).Value = "Inactive"
$ws.Cells.Item(320, 1).Value = "2026-02-04"
$ws.Cells.Item(320, 2).Value = "14:27:43"
$ws.Cells.Item(320, 3).Value = "14:00"
$ws.Cells.Item(320, 4).Value = "Bathroom"
$ws.Cells.Item(320, 5).Value = "No Motion"
$ws.Cells.Item(320, 6).Value = "Inactive"
$ws.Cells.Item(321, 1).Value = "2026-02-04"
$ws.Cells.Item(321, 2).Value = "14:27:48"
$ws.Cells.Item(321, 3).Value = "14:00"
$ws.Cells.Item(321, 4).Value = "Bathroom"
$ws.Cells.Item(321, 5).Value = "No Motion"
$ws.Cells.Item(321, 6).Value = "Inactive"
$ws.Cells.Item(322, 1).Value = "2026-02-04"
$ws.Cells.Item(322, 2).Value = "14:27:53"
$ws.Cells.Item(322, 3).Value = "14:00"
$ws.Cells.Item(322, 4).Value = "Bathroom"
$ws.Cells.Item(322, 5).Value = "No Motion"
$ws.Cells.Item(322, 6).Value = "Inactive"
$ws.Cells.Item(323, 1).Value = "2026-02-04"
$ws.Cells.Item(323, 2).Value = "14:27:58"
$ws.Cells.Item(323, 3).Value = "14:00"
$ws.Cells.Item(323, 4).Value = "Bathroom"
$ws.Cells.Item(323, 5).Value = "No Motion"
$ws.Cells.Item(323, 6).Value = "Inactive"
$ws.Range("A312:F323").ClearFormats()

$ws = $wb.Worksheets.Item("Humidity")
$ws.Range("A255:F265").NumberFormat = "@"
$ws.Cells.Item(255, 1).Value = "2026-02-04"
$ws.Cells.Item(255, 2).Value = "14:27:01"
$ws.Cells.Item(255, 3).Value = "14:00"
$ws.Cells.Item(255, 4).Value = "Bathroom"
$ws.Cells.Item(255, 5).Value = "79.4%"
$ws.Cells.Item(255, 6).Value = "Active"
$ws.Cells.Item(256, 1).Value = "2026-02-04"
$ws.Cells.Item(256, 2).Value = "14:27:06"
$ws.Cells.Item(256, 3).Value = "14:00"
$ws.Cells.Item(256, 4).Value = "Bathroom"
$ws.Cells.Item(256, 5).Value = "78.5%"
$ws.Cells.Item(256, 6).Value = "Active"
$ws.Cells.Item(257, 1).Value = "2026-02-04"
$ws.Cells.Item(257, 2).Value = "14:27:16"
$ws.Cells.Item(257, 3).Value = "14:00"
$ws.Cells.Item(257, 4).Value = "Bathroom"
$ws.Cells.Item(257, 5).Value = "78.5%"
$ws.Cells.Item(257, 6).Value = "Active"
$ws.Cells.Item(258, 1).Value = "2026-02-04"
$ws.Cells.Item(258, 2).Value = "14:27:21"
$ws.Cells.Item(258, 3).Value = "14:00"
$ws.Cells.Item(258, 4).Value = "Bathroom"
$ws.Cells.Item(258, 5).Value = "79.4%"
$ws.Cells.Item(258, 6).Value = "Active"
$ws.Cells.Item(259, 1).Value = "2026-02-04"
$ws.Cells.Item(259, 2).Value = "14:27:26"
$ws.Cells.Item(259, 3).Value = "14:00"
$ws.Cells.Item(259, 4).Value = "Bathroom"
$ws.Cells.Item(259, 5).Value = "78.5%"
$ws.Cells.Item(259, 6).Value = "Active"
$ws.Cells.Item(260, 1).Value = "2026-02-04"
$ws.Cells.Item(260, 2).Value = "14:27:31"
$ws.Cells.Item(260, 3).Value = "14:00"
$ws.Cells.Item(260, 4).Value = "Bathroom"
$ws.Cells.Item(260, 5).Value = "79.5%"
$ws.Cells.Item(260, 6).Value = "Active"
$ws.Cells.Item(261, 1).Value = "2026-02-04"
$ws.Cells.Item(261, 2).Value = "14:27:36"
$ws.Cells.Item(261, 3).Value = "14:00"
$ws.Cells.Item(261, 4).Value = "Bathroom"
$ws.Cells.Item(261, 5).Value = "78.5%"
$ws.Cells.Item(261, 6).Value = "Active"
$ws.Cells.Item(262, 1).Value = "2026-02-04"
$ws.Cells.Item(262, 2).Value = "14:27:41"
$ws.Cells.Item(262, 3).Value = "14:00"
$ws.Cells.Item(262, 4).Value = "Bathroom"
$ws.Cells.Item(262, 5).Value = "79.4%"
$ws.Cells.Item(262, 6).Value = "Active"
$ws.Cells.Item(263, 1).Value = "2026-02-04"
$ws.Cells.Item(263, 2).Value = "14:27:46"
$ws.Cells.Item(263, 3).Value = "14:00"
$ws.Cells.Item(263, 4).Value = "Bathroom"
$ws.Cells.Item(263, 5).Value = "78.5%"
$ws.Cells.Item(263, 6).Value = "Active"
$ws.Cells.Item(264, 1).Value = "2026-02-04"
$ws.Cells.Item(264, 2).Value = "14:27:51"
$ws.Cells.Item(264, 3).Value = "14:00"
$ws.Cells.Item(264, 4).Value = "Bathroom"
$ws.Cells.Item(264, 5).Value = "79.4%"
$ws.Cells.Item(264, 6).Value = "Active"
$ws.Cells.Item(265, 1).Value = "2026-02-04"
$ws.Cells.Item(265, 2).Value = "14:27:56"
$ws.Cells.Item(265, 3).Value = "14:00"
$ws.Cells.Item(265, 4).Value = "Bathroom"
$ws.Cells.Item(265, 5).Value = "78.5%"
$ws.Cells.Item(265, 6).Value = "Active"
$ws.Range("A255:F265").ClearFormats()

$ws = $wb.Worksheets.Item("Temperature")
$ws.Range("A255:F265").NumberFormat = "@"
$ws.Cells.Item(255, 1).Value = "2026-02-04"
$ws.Cells.Item(255, 2).Value = "14:27:02"
$ws.Cells.Item(255, 3).Value = "14:00"
$ws.Cells.Item(255, 4).Value = "Bathroom"
$ws.Cells.Item(255, 5).Value = "24.3C"
$ws.Cells.Item(255, 6).Value = "Active"
$ws.Cells.Item(256, 1).Value = "2026-02-04"
$ws.Cells.Item(256, 2).Value = "14:27:06"
$ws.Cells.Item(256, 3).Value = "14:00"
$ws.Cells.Item(256, 4).Value = "Bathroom"
$ws.Cells.Item(256, 5).Value = "24.3C"
$ws.Cells.Item(256, 6).Value = "Active"
$ws.Cells.Item(257, 1).Value = "2026-02-04"
$ws.Cells.Item(257, 2).Value = "14:27:17"
$ws.Cells.Item(257, 3).Value = "14:00"
$ws.Cells.Item(257, 4).Value = "Bathroom"
$ws.Cells.Item(257, 5).Value = "24.3C"
$ws.Cells.Item(257, 6).Value = "Active"
$ws.Cells.Item(258, 1).Value = "2026-02-04"
$ws.Cells.Item(258, 2).Value = "14:27:22"
$ws.Cells.Item(258, 3).Value = "14:00"
$ws.Cells.Item(258, 4).Value = "Bathroom"
$ws.Cells.Item(258, 5).Value = "24.3C"
$ws.Cells.Item(258, 6).Value = "Active"
$ws.Cells.Item(259, 1).Value = "2026-02-04"
$ws.Cells.Item(259, 2).Value = "14:27:27"
$ws.Cells.Item(259, 3).Value = "14:00"
$ws.Cells.Item(259, 4).Value = "Bathroom"
$ws.Cells.Item(259, 5).Value = "24.3C"
$ws.Cells.Item(259, 6).Value = "Active"
$ws.Cells.Item(260, 1).Value = "2026-02-04"
$ws.Cells.Item(260, 2).Value = "14:27:32"
$ws.Cells.Item(260, 3).Value = "14:00"
$ws.Cells.Item(260, 4).Value = "Bathroom"
$ws.Cells.Item(260, 5).Value = "24.3C"
$ws.Cells.Item(260, 6).Value = "Active"
$ws.Cells.Item(261, 1).Value = "2026-02-04"
$ws.Cells.Item(261, 2).Value = "14:27:37"
$ws.Cells.Item(261, 3).Value = "14:00"
$ws.Cells.Item(261, 4).Value = "Bathroom"
$ws.Cells.Item(261, 5).Value = "24.3C"
$ws.Cells.Item(261, 6).Value = "Active"
$ws.Cells.Item(262, 1).Value = "2026-02-04"
$ws.Cells.Item(262, 2).Value = "14:27:42"
$ws.Cells.Item(262, 3).Value = "14:00"
$ws.Cells.Item(262, 4).Value = "Bathroom"
$ws.Cells.Item(262, 5).Value = "24.3C"
$ws.Cells.Item(262, 6).Value = "Active"
$ws.Cells.Item(263, 1).Value = "2026-02-04"
$ws.Cells.Item(263, 2).Value = "14:27:47"
$ws.Cells.Item(263, 3).Value = "14:00"
$ws.Cells.Item(263, 4).Value = "Bathroom"
$ws.Cells.Item(263, 5).Value = "24.4C"
$ws.Cells.Item(263, 6).Value = "Active"
$ws.Cells.Item(264, 1).Value = "2026-02-04"
$ws.Cells.Item(264, 2).Value = "14:27:52"
$ws.Cells.Item(264, 3).Value = "14:00"
$ws.Cells.Item(264, 4).Value = "Bathroom"
$ws.Cells.Item(264, 5).Value = "24.3C"
$ws.Cells.Item(264, 6).Value = "Active"
$ws.Cells.Item(265, 1).Value = "2026-02-04"
$ws.Cells.Item(265, 2).Value = "14:27:57"
$ws.Cells.Item(265, 3).Value = "14:00"
$ws.Cells.Item(265, 4).Value = "Bathroom"
$ws.Cells.Item(265, 5).Value = "24.4C"
$ws.Cells.Item(265, 6).Value = "Active"
$ws.Range("A255:F265").ClearFormats()

